$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.716.92"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").Value = "'1.805.40"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'231.94"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").Value = "'0.5947"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.2787"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "'0.06844"
$ws.Range("E9").Value = "  -3.50%  "
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").Value = "'0.07525"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").Value = "'1.793.16"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").Value = "'4.814"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").Value = "'0.6241"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "'2.050.68"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("E16").Value = "  -6.66%  "
$ws.Range("D17").Value = "'75.81"
$ws.Range("E17").Value = "  -3.56%  "
$ws.Range("D18").Value = "'28.666.35"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").Value = "'5.507"
$ws.Range("E19").Value = "  -6.06%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'210.90"
$ws.Range("E21").Value = "  -6.96%  "
$ws.Range("D22").Value = "'11.48"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("D23").Value = "'6.874"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "'154.32"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").Value = "'7.884"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").Value = "'0.1276"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("D28").Value = "'16.44"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("E29").Value = "  -4.45%  "
$ws.Range("D30").Value = "'0.06247"
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("D31").Value = "'1.423"
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("D32").Value = "'3.788"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "'3.763"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").Value = "'1.724"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").Value = "'1.068"
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("D36").Value = "'0.6426"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").Value = "'2.494"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").Value = "'6.496"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").Value = "'1.140.15"
$ws.Range("E41").Value = "  -6.10%  "
$ws.Range("D42").Value = "'0.8794"
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("D43").Value = "'1.008"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").Value = "'100.67"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "'1.966.99"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("E47").Value = "  -4.71%  "
$ws.Range("D48").Value = "'1.607"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "'8.353"
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").Value = "'0.05464"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").Value = "'0.4484"
$ws.Range("E51").Value = "  -1.85%  "
